$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet 1: "ملخص المشروع" (project summary) — insert a new SSOT line
# "ceiling ceramic" above "إجمالي السيراميك (SSOT)", rename/adjust a
# few of the rows that follow, and add the missing unit for the
# "opening count" row.
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Push the existing rows 11-17 down to 12-18, leaving a blank row 11
# (this also copies the B2/C3 style formatting down automatically).
$ws1.Rows.Item(11).Insert()

# New row 11: "سيراميك الأسقف (SSOT)"
$ws1.Cells.Item(11, 2).Value = "سيراميك الأسقف (SSOT)"
$ws1.Cells.Item(11, 3).Value = 0
$ws1.Cells.Item(11, 4).Value = "م²"

# Row 13 (was "نعلات (SSOT)" at old row 12): value resets to 0
$ws1.Cells.Item(13, 3).Value = 0

# Row 14 (was "حجر/أطر (SSOT)" at old row 13): rename label
$ws1.Cells.Item(14, 2).Value = "حجر/أطر (مجموع الغرف)"

# Row 16 (was "عدد الفتحات (أبواب+شبابيك)" at old row 15): add unit
$ws1.Cells.Item(16, 4).Value = "قطعة"

# ------------------------------------------------------------------
# Sheet 2: "مساحة الغرف" (room areas) — add a new "ملاحظات النعلات"
# column (S) and zero-out the "نعلات" quantity for room 1 since it is
# now cancelled out (wall ceramic exists instead).
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# New column S mirrors the width + header/body formatting of the
# existing last column (R).
$ws2.Columns.Item(19).ColumnWidth = $ws2.Columns.Item(18).ColumnWidth

$ws2.Cells.Item(1, 18).Copy($ws2.Cells.Item(1, 19))
$ws2.Cells.Item(1, 19).Value = "ملاحظات النعلات"

$ws2.Cells.Item(2, 2).Copy($ws2.Cells.Item(2, 19))
$ws2.Cells.Item(2, 19).Value = "ملغى (يوجد سيراميك جدران)"

# "نعلات" quantity for room 1 drops to 0 (cancelled, see note above)
$ws2.Cells.Item(2, 17).Value = 0
